$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells whose new value is a "pure number" string (e.g. "244.08", "0.630")
# must be forced to Text, otherwise Excel/COM auto-converts the assignment
# to a numeric value and mangles formatting (drops trailing zeros, turns
# "231.00" into 231, adds floating point noise, etc). We force text with a
# leading apostrophe (the standard Excel 'text qualifier' trick) and then
# reset the cell style back to Normal so no stray quotePrefix style sticks.

$ws.Range('D2').Value = '42.187.81'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '2.228.60'
$ws.Range('E3').Value = '  -0.77%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = "'244.08"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.12%  '
$ws.Range('D6').Value = "'0.630"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.28%  '
$ws.Range('D7').Value = "'73.83"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.68%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').Value = "'0.617"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('D10').Value = "'43.32"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.68%  '
$ws.Range('D11').Value = "'0.0966"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.63%  '
$ws.Range('D12').Value = "'7.16"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.11%  '
$ws.Range('E13').Value = '  +1.25%  '
$ws.Range('D14').Value = "'14.32"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.15%  '
$ws.Range('D15').Value = "'0.848"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('D16').Value = '2.226.42'
$ws.Range('E16').Value = '  -0.46%  '
$ws.Range('D17').Value = '42.071.86'
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('E18').Value = '  +14.35%  '
$ws.Range('D19').Value = "'6.24"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.46%  '
$ws.Range('D20').Value = "'72.45"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.08%  '
$ws.Range('D21').Value = "'10.37"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +38.92%  '
$ws.Range('D22').Value = "'231.00"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('E23').Value = '  -7.89%  '
$ws.Range('D24').Value = "'11.76"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.95%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('E26').Value = '  +1.22%  '
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('E28').Value = '  +3.03%  '
$ws.Range('D29').Value = "'166.96"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.54%  '
$ws.Range('D30').Value = "'20.71"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = "'5.63"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +16.38%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = "'0.0802"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.78%  '
$ws.Range('E33').Value = '  -3.45%  '
$ws.Range('D34').Value = "'0.125"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.38%  '
$ws.Range('D35').Value = "'29.56"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.72%  '
$ws.Range('D36').Value = "'4.38"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.64%  '
$ws.Range('D37').Value = "'0.0304"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.13%  '
$ws.Range('D38').Value = "'13.18"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.14%  '
$ws.Range('D39').Value = "'2.16"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.46%  '
$ws.Range('D40').Value = "'64.90"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.88%  '
$ws.Range('D41').Value = "'5.60"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.95%  '
$ws.Range('D42').Value = "'0.202"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.36%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = "'105.49"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.77%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = "'8.75"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.29%  '
$ws.Range('D46').Value = "'2.43"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +8.11%  '
$ws.Range('D47').Value = "'1.13"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.44%  '
$ws.Range('D48').Value = "'1.17"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.34%  '
$ws.Range('E49').Value = '  +0.87%  '
$ws.Range('B50').Value = 'SynthetixNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D50').Value = "'4.07"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.436.32'
$ws.Range('E51').Value = '  -0.64%  '
